$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not be coerced to a
# number by Excel's auto-detection) while keeping its original style (s="8").
# Use a leading apostrophe to force text, then re-apply the original cell
# formatting (copied from the neighboring text cell B2) so the style index
# does not drift to a newly generated "quote prefixed" style.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null

$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 15.09.2024"

# Transaction row 6
$ws.Range("B6").Value = "19.09."
$ws.Range("C6").Value = "20.09."
$ws.Range("D6").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E6").Value = "85,07-"

# Transaction row 7
$ws.Range("B7").Value = "22.09."
$ws.Range("C7").Value = "23.09."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,23-"

# Transaction row 8
$ws.Range("B8").Value = "26.09."
$ws.Range("C8").Value = "27.09."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 91671471"
$ws.Range("E8").Value = "42,18-"

# Transaction row 9 - was previously empty, now a new transaction.
# Copy formatting from row 8 (the last populated transaction row) before
# setting the values so the new cells pick up the correct style (s="17"
# on E9, matching the other amount cells, instead of the empty-row style).
$ws.Range("B8:E8").Copy() | Out-Null
$ws.Range("B9:E9").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Value = "30.09."
$ws.Range("C9").Value = "01.10."
$ws.Range("D9").Value = "KARTENZ./30.09 REWE RO"
$ws.Range("E9").Value = "69,21-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 05.10.2024"
$ws.Range("E12").Value = "221,69-"

# Next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.10.2024"
